# Generate Report for Handoff
# Updates the localization-status workbook to reflect that b.md has been
# handed off (new xliff files generated) while a.md's status text is
# refreshed to the current wording.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1355de321760fc871c649d70891e8cf9dc175ced/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e4ca56aa2f27370d102f6763d4c59070d2479ea1/e2e/b.md."

# ---------------------------------------------------------------------
# Overview sheet - b.md row (row 3): status + latest generate date
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-09-05 00:43:10"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

# a.md (row 2) - status text refreshed
$zhcn.Range("C2").Value = "Ready for handoff"

# b.md (row 3) - new handoff file generated
$zhcn.Range("C3").Value = "Ready for handoff"
# Leading apostrophe forces text (otherwise "False" is auto-coerced to a
# boolean cell); reapplying the Normal style clears the resulting
# quote-prefix formatting flag so the cell round-trips as plain text.
$zhcn.Range("F3").Value = "'False"
$zhcn.Range("F3").Style = "Normal"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-09-05 00:43:01"
$zhcn.Range("P3").Value = $errorDetail

# Widen the Error Detail column to fit the new message
$zhcn.Columns.Item(16).ColumnWidth = 39.1666666666667

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

# a.md (row 2) - status text + generate date refreshed
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-09-05 00:43:10"

# b.md (row 3) - new handoff file generated
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Value = "'False"
$dede.Range("F3").Style = "Normal"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-09-05 00:43:10"
$dede.Range("P3").Value = $errorDetail

# Widen the Error Detail column to fit the new message
$dede.Columns.Item(16).ColumnWidth = 39.1666666666667
